# Update IBAN/bank-account test data values and adjust the active selection
# as recorded in the workbook's sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "CH9689144715152235363"
$ws.Range("B15").Value = "DE03500105177284191473"
$ws.Range("B16").Value = "DE64500105177324151368"
$ws.Range("B17").Value = "IT37T0300203280682244199423"

# Move the active cell / selection to B17, matching the saved sheet view.
$ws.Range("B17").Select()
